$wb = $excel.ActiveWorkbook

$wsReq   = $wb.Worksheets.Item(1)   # "Требования"
$wsCases = $wb.Worksheets.Item(2)   # "Тест-кейсы"

# ---------------------------------------------------------------------------
# 1. "Требования" sheet: a new requirement "F-18" is inserted, which shifts
#    the existing F-14 / F-15 / F-16 / F-17 identifiers in rows 37-40 down by
#    one, and the last two rows (39 & 40) pick up the "full border" look
#    already used by the rows above them (A34 is a cell that already has
#    that exact style).
# ---------------------------------------------------------------------------
$wsReq.Cells.Item(37, 1).Value = "F-15"
$wsReq.Cells.Item(38, 1).Value = "F-16"
$wsReq.Cells.Item(39, 1).Value = "F-17"
$wsReq.Cells.Item(40, 1).Value = "F-18"

$fmtSrc = $wsReq.Cells.Item(34, 1)
$fmtSrc.Copy()
$wsReq.Cells.Item(39, 1).PasteSpecial(-4122)
$wsReq.Cells.Item(40, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. "Тест-кейсы" sheet: column H (requirement reference) gets populated for
#    rows 5-41, and A40 is renamed from "RF-13" to the new "RF-14".
# ---------------------------------------------------------------------------
$wsCases.Cells.Item(5, 8).Value  = "T-5"
$wsCases.Cells.Item(6, 8).Value  = "T-5"
$wsCases.Cells.Item(7, 8).Value  = "T-6"
$wsCases.Cells.Item(8, 8).Value  = "T-7"
$wsCases.Cells.Item(9, 8).Value  = "T-8"
$wsCases.Cells.Item(10, 8).Value = "T-9"
$wsCases.Cells.Item(11, 8).Value = "T-10"

$wsCases.Cells.Item(12, 8).Value = "F-3/1"
$wsCases.Cells.Item(13, 8).Value = "F-3/2"
$wsCases.Cells.Item(14, 8).Value = "F-3/3"
$wsCases.Cells.Item(15, 8).Value = "F-3/4"

$wsCases.Cells.Item(16, 8).Value = "F-4"
$wsCases.Cells.Item(17, 8).Value = "F-4"
$wsCases.Cells.Item(18, 8).Value = "F-4"
$wsCases.Cells.Item(19, 8).Value = "F-4"

$wsCases.Cells.Item(20, 8).Value = "F-5"
$wsCases.Cells.Item(21, 8).Value = "F-5"
$wsCases.Cells.Item(22, 8).Value = "F-5"
$wsCases.Cells.Item(23, 8).Value = "F-5"

$wsCases.Cells.Item(24, 8).Value = "F-6"
$wsCases.Cells.Item(25, 8).Value = "F-6"
$wsCases.Cells.Item(26, 8).Value = "F-6"
$wsCases.Cells.Item(27, 8).Value = "F-6"

$wsCases.Cells.Item(28, 8).Value = "F-7"
$wsCases.Cells.Item(29, 8).Value = "F-7"
$wsCases.Cells.Item(30, 8).Value = "F-7"
$wsCases.Cells.Item(31, 8).Value = "F-7"

$wsCases.Cells.Item(32, 8).Value = "F-8"
$wsCases.Cells.Item(33, 8).Value = "F-8"
$wsCases.Cells.Item(34, 8).Value = "F-8"
$wsCases.Cells.Item(35, 8).Value = "F-8"

$wsCases.Cells.Item(36, 8).Value = "F-9"
$wsCases.Cells.Item(37, 8).Value = "F-10"
$wsCases.Cells.Item(38, 8).Value = "F-12"
$wsCases.Cells.Item(39, 8).Value = "F-13"

$wsCases.Cells.Item(40, 1).Value = "RF-14"
$wsCases.Cells.Item(40, 8).Value = "F-14"

$wsCases.Cells.Item(41, 8).Value = "F-18"

# ---------------------------------------------------------------------------
# 3. View state: update the remembered selection on each sheet, then leave
#    "Тест-кейсы" as the active (displayed) sheet/tab, matching the saved
#    workbook view.
# ---------------------------------------------------------------------------
$wsReq.Activate()
$wsReq.Range("B48").Select()

$wsCases.Activate()
$wsCases.Range("D39").Select()

Write-Host "done"
